$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParaXml($para, [string]$innerAttrsAndContent) {
    $xml = "<w:p $wNs $innerAttrsAndContent</w:p>"
    $result = $para.Range.InsertXML($xml)
}

# Locate the target paragraphs by distinctive, stable text content.
$paraAstor = $null
$paraPydantic = $null
$paraNiUiCreator = $null
$paraDevVersion = $null
$paraBooleans = $null
$paraConverterCmd = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "astor*") { $paraAstor = $p }
    elseif ($t -like "pydantic*") { $paraPydantic = $p }
    elseif ($t -like "ni-measurement-ui-creator*") { $paraNiUiCreator = $p }
    elseif ($t -like "1.0.0.dev*or above*") { $paraDevVersion = $p }
    elseif ($t -like "*doesn*support List of strings and List of booleans*") { $paraBooleans = $p }
    elseif ($t -like "ni-measurement-plugin-converter --display-name*") { $paraConverterCmd = $p }
}

# 1) "astor" table cell - drop the spellStart/spellEnd proofErr wrap.
Set-ParaXml $paraAstor 'w14:paraId="5C5D0AC9" w14:textId="7FF1ABBC" w:rsidR="00A63694" w:rsidRDefault="00A63694" w:rsidP="001E1678"><w:r><w:t>astor</w:t></w:r>'

# 2) "pydantic" table cell - drop the spellStart/spellEnd proofErr wrap.
Set-ParaXml $paraPydantic 'w14:paraId="2786EA56" w14:textId="16D65B87" w:rsidR="00A63694" w:rsidRDefault="00A63694" w:rsidP="001E1678"><w:r><w:t>pydantic</w:t></w:r>'

# 3) "ni-measurement-ui-creator" table cell - drop both proofErr wraps (around "n"+"i" and around "ui").
Set-ParaXml $paraNiUiCreator 'w14:paraId="12C12DF9" w14:textId="36852335" w:rsidR="00A63694" w:rsidRDefault="009A1BE3" w:rsidP="00A63694"><w:r><w:t>n</w:t></w:r><w:r w:rsidR="00A63694"><w:t>i</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r w:rsidR="00A63694"><w:t>measurement</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r w:rsidR="00A63694"><w:t>ui</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r w:rsidR="00A63694"><w:t>creator</w:t></w:r>'

# 4) "1.0.0.dev5 or above" -> "1.0.0.dev6 or above"
Set-ParaXml $paraDevVersion 'w14:paraId="180D6015" w14:textId="2EE260F4" w:rsidR="00A63694" w:rsidRDefault="00A63694" w:rsidP="001E1678"><w:r><w:t>1.0.0.dev</w:t></w:r><w:r w:rsidR="00E52AA1"><w:t>6</w:t></w:r><w:r><w:t xml:space="preserve"> or above</w:t></w:r>'

# 5) Booleans bullet - merge the three runs (text / "booleans" / ".") into one run and drop proofErr.
$rsquo = [char]0x2019
$booleansXml = 'w14:paraId="720CCE32" w14:textId="46421D73" w:rsidR="001E1678" w:rsidRDefault="001E1678" w:rsidP="004B1380"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:t>Measurement UI created by the tool doesn' + $rsquo + 't support List of strings and List of booleans.</w:t></w:r>'
Set-ParaXml $paraBooleans $booleansXml

# 6) "ni-measurement-plugin-converter ..." bullet - drop the "n"/"i" proofErr wrap and merge the
#    "dir"-split runs for --measurement-file-dir / --output-dir into their neighbours.
Set-ParaXml $paraConverterCmd 'w14:paraId="0CEAC792" w14:textId="6E0E844A" w:rsidR="00A42C2F" w:rsidRDefault="009A1BE3" w:rsidP="007A098F"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>n</w:t></w:r><w:r w:rsidR="008006C1"><w:t>i</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r w:rsidR="008006C1"><w:t>measurement</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r w:rsidR="008006C1"><w:t>plugin</w:t></w:r><w:r><w:t>-</w:t></w:r><w:r w:rsidR="008006C1"><w:t xml:space="preserve">converter --display-name &lt;Display Name&gt; --measurement-file-dir &lt;Measurement file directory&gt; --function &lt;Measurement function name&gt; </w:t></w:r><w:r w:rsidR="00553E9C"><w:t>--output-dir &lt;Output directory&gt;</w:t></w:r>'

Write-Output "done"
